# Build site at 2023-01-09 16:18:13 UTC
# - "Ativacao:" date updated from 01/01/2018 to 01/01/2023
# - "Semestral" entry removed (its slot now reuses the updated date text)
# - "Metodo:" teacher entry updated from "5840712 - Angelo Capri Neto"
#   to the new "5840521 - Rosa Ana Conte" (old teacher value now also
#   appears where the removed "Semestral" shared-string shift lands)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellAddress,
        [string]$Text
    )
    # Write through a formula (so Excel does not try to auto-detect the
    # text as a date/number) and then Copy/PasteSpecial values-only onto
    # the real target cell, so the destination's existing style/number
    # format is left completely untouched.
    $scratch = $ws.Range("ZZ1")
    $escaped = $Text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($CellAddress).PasteSpecial(-4163) # xlPasteValues
    $scratch.ClearContents()
}

# Ativacao: 01/01/2018 -> 01/01/2023
Set-TextValue "B8" "01/01/2023"
Set-TextValue "C8" "01/01/2023"

# Row that held "Semestral" now carries the (new) activation date text
Set-TextValue "B13" "01/01/2023"
Set-TextValue "C13" "01/01/2023"

# Row that held the old activation date now carries the original teacher
Set-TextValue "B15" "5840712 - Ângelo Capri Neto"
Set-TextValue "C15" "5840712 - Ângelo Capri Neto"

# Metodo: responsible teacher changes to the new entry
Set-TextValue "B18" "5840521 - Rosa Ana Conte"
Set-TextValue "C18" "5840521 - Rosa Ana Conte"

$excel.DisplayAlerts = $false
$wb.Save()
